# Generate Report for Handoff
# The file "8e8a11af-c2bc-445c-a82c-9a46df9ff85d.md" moved from
# "Handed back: in sync with en-US" to "Ready for handoff" status, with a
# fresh handoff timestamp recorded per-locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 8e8a11af-...md (row 3), zh-cn (B3) & de-de (C3) columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for 8e8a11af-...md (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-25 06:31:49"

# --- de-de sheet: row for 8e8a11af-...md (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-25 06:32:00"
